# Fruta / hortaliza, semanal
# Insert two new weekly price rows (week ending 44516) for "Feria Lagunitas de
# Puerto Montt" - Naranja, pushing the existing rows 266-288 down to 268-290.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 266 (existing data shifts down to 268..290).
$ws.Range("A266:A267").EntireRow.Insert()

# New row 266: Naranja, Lane Late, Primera - Provincia de Melipilla
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44516
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = "Fruta"
$ws.Range("G266").Value = 100102
$ws.Range("H266").Value = "Cítricos"
$ws.Range("I266").Value = 100102005
$ws.Range("J266").Value = "Naranja"
$ws.Range("K266").Value = "Lane Late"
$ws.Range("L266").Value = "Primera"
$ws.Range("M266").Value = 600
$ws.Range("N266").Value = 12000
$ws.Range("O266").Value = 13000
$ws.Range("P266").Value = 12500
$ws.Range("Q266").Value = "$/caja 15 kilos empedrada"
$ws.Range("R266").Value = "Provincia de Melipilla"
$ws.Range("S266").Value = 833
$ws.Range("T266").Value = 15

# New row 267: Naranja, Navel Late, Primera - Provincia de Melipilla
$ws.Range("A267").Value = 4
$ws.Range("B267").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C267").Value = "Los Lagos"
$ws.Range("D267").Value = 44516
$ws.Range("E267").Value = 10
$ws.Range("F267").Value = "Fruta"
$ws.Range("G267").Value = 100102
$ws.Range("H267").Value = "Cítricos"
$ws.Range("I267").Value = 100102005
$ws.Range("J267").Value = "Naranja"
$ws.Range("K267").Value = "Navel Late"
$ws.Range("L267").Value = "Primera"
$ws.Range("M267").Value = 600
$ws.Range("N267").Value = 12000
$ws.Range("O267").Value = 13000
$ws.Range("P267").Value = 12500
$ws.Range("Q267").Value = "$/caja 15 kilos empedrada"
$ws.Range("R267").Value = "Provincia de Melipilla"
$ws.Range("S267").Value = 833
$ws.Range("T267").Value = 15
